# Updated Global_M2 for easier usage.
# Update existing row 308 (D and F changed), and append new rows 309-311
# with the same FX_IDC:USDILS OHLC data layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 308 ---
$ws.Cells.Item(308, 4).Value = 3.6919
$ws.Cells.Item(308, 6).Value = 3.6137

# --- New rows data: row, datetime(serial), symbol, open, high, low, close, volume ---
$newRows = @(
    @(309, 45047.33333333334, "FX_IDC:USDILS", 3.624,  3.74944, 3.6065, 3.7154, 0),
    @(310, 45078.33333333334, "FX_IDC:USDILS", 3.7155, 3.7794,  3.5401, 3.6932, 0),
    @(311, 45110.33333333334, "FX_IDC:USDILS", 3.7059, 3.72694, 3.6858, 3.6949, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Copy formatting (incl. date style) from the row above into the new row
    $ws.Range("A" + ($r - 1) + ":G" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
